$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade row (row 8) mirroring the existing trade rows (3-7).
$row = 8

$ws.Cells.Item($row, 1).Value = 42636.592766203707
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = $true
$ws.Cells.Item($row, 3).Value = 9860.5
$ws.Cells.Item($row, 4).Value = 9830.52
$ws.Cells.Item($row, 5).Value = 81.97
$ws.Cells.Item($row, 6).Value = 81.47

$ws.Cells.Item($row, 7).Value = $true
$ws.Cells.Item($row, 7).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 8).Value = -0.61
$ws.Cells.Item($row, 9).Value = $false
